# Apply the changes described by the diff:
#  - sheet3 "Blocks per Program": drop the old column-B "Total" (SUM) column,
#    fix the "?" placeholders in I2/J2/I3/J3 with real numbers, and add a new
#    column-K "Total" (SUM) column instead; update the sheet's selection.
#  - sheet4 "Sheet1": append a new transposed summary table (rows 13-21) and
#    make this sheet the active one.
#  - workbook: active tab moves from "Blocks per Program" to "Sheet1".

$wb  = $excel.ActiveWorkbook
$wsBlocks = $wb.Worksheets.Item("Blocks per Program")
$wsSheet1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. sheet3 ("Blocks per Program")
# ---------------------------------------------------------------------------

# Fix the four cells that used to hold the "?" placeholder text.
$wsBlocks.Range("I2").Value = 0
$wsBlocks.Range("J2").Value = 3
$wsBlocks.Range("I3").Value = 0
$wsBlocks.Range("J3").Value = 0

# Drop the old "Total" column (B) - header + all per-row SUM formulas.
$wsBlocks.Range("B1:B18").ClearContents() | Out-Null

# Add the new "Total" column (K) - same header text, same SUM formula shape,
# just relocated after the data columns (C:J) instead of before them.
$wsBlocks.Range("K1").Value = "Total"
for ($r = 2; $r -le 18; $r++) {
    $wsBlocks.Range("K$r").Formula = "=SUM(C$r`:J$r)"
}

# ---------------------------------------------------------------------------
# 2. sheet4 ("Sheet1") - new rows 13-21
# ---------------------------------------------------------------------------

# Row label (col A) -> category name already present in sharedStrings.
$rowLabels = @{
    13 = "Data"
    14 = "Actor"
    15 = "Sensor"
    16 = "Logic"
    17 = "MyBlock Call"
    18 = "Comment"
    19 = "Variables"
    20 = "MyBlocks"
}

# data[row][col] : col 2..18 -> B..R, one column per program (transposed
# from sheet3's per-program rows 2-18, columns C:J).
$data = @{
    13 = @(1,0,2,5,0,0,4,4,1,0,0,2,0,7,1,1,0)
    14 = @(8,6,5,7,16,8,10,16,10,1,2,1,7,15,2,2,78)
    15 = @(1,0,0,2,1,0,8,10,4,1,0,1,0,4,2,0,0)
    16 = @(3,5,4,4,4,3,7,11,2,2,1,4,11,2,0,0,25)
    17 = @(2,0,2,2,3,4,3,5,0,1,0,1,0,3,1,4,27)
    18 = @(2,0,0,1,0,0,1,0,0,6,6,0,15,0,0,7,13)
    19 = @(0,0,1,3,1,1,0,0,0,0,0,0,0,19,0,0,0)
    20 = @(0,3,0,2,3,2,2,4,0,1,0,1,0,4,1,3,6)
}

foreach ($r in 13..20) {
    $wsSheet1.Range("A$r").Value = $rowLabels[$r]
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $col = $i + 2   # column B = 2
        $wsSheet1.Cells.Item($r, $col).Value = $rowVals[$i]
    }
}

# Row 21: "Total" row, column-wise SUM of rows 13-20.
$wsSheet1.Range("A21").Value = "Total"
foreach ($col in 2..18) {
    $colLetter = [char](64 + $col)
    $wsSheet1.Cells.Item(21, $col).Formula = "=SUM($colLetter`13:$colLetter`20)"
}

# ---------------------------------------------------------------------------
# 3. Selections + active sheet/tab
# ---------------------------------------------------------------------------

# Update sheet3's selection (this sheet is no longer the active one).
$wsBlocks.Range("C1:K18").Select() | Out-Null

# Sheet1 becomes the active sheet/tab, with its own new selection.
$wsSheet1.Activate() | Out-Null
$wsSheet1.Range("A13:R21").Select() | Out-Null
